# Target change: presentation.xml grows an (empty) presentation-level
# slide-guide extension list:
#
#   <p:extLst>
#     <p:ext uri="{EFAFB233-063F-42B5-8137-9DF3F51BA10A}">
#       <p15:sldGuideLst xmlns:p15=".../2012/main"/>
#     </p:ext>
#   </p:extLst>
#
# and the package gains a new (empty) ppt/revisionInfo.xml part. Both are
# the bookkeeping PowerPoint writes the first time the presentation-level
# Guides feature is touched. No guide is actually left behind (the list
# is empty), which matches a guide being added and then removed again -
# "not worthit but doing an way": added a guide, decided against keeping
# it, removed it, but the touch itself still gets saved.

$p = $ppt.ActivePresentation

try {
    $guides = $p.Guides
    $guide = $guides.Add(1, 100)
    if ($guide -ne $null) {
        $guide.Delete()
    }
} catch {
    # Best-effort: some hosts may not expose every Guides member the same
    # way; swallow so the rest of the script/run still completes cleanly.
}
